$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text while we overwrite them,
# matching the workbook author's original inlineStr cells (values like "0.999" or
# "1.00" would otherwise be auto-coerced to numbers and lose their literal formatting).
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.446.49"
$ws.Range("E2").Value = "  -4.06%  "
$ws.Range("D3").Value = "3.309.21"
$ws.Range("E3").Value = "  -4.19%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "558.88"
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("D6").Value = "143.43"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.308.88"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").Value = "0.410"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "3.869.09"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D15").Value = "27.25"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "3.291.91"
$ws.Range("E16").Value = "  -4.83%  "
$ws.Range("E17").Value = "  -4.13%  "
$ws.Range("D18").Value = "60.355.63"
$ws.Range("E18").Value = "  -4.14%  "
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  -5.06%  "
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "8.59"
$ws.Range("E21").Value = "  -4.97%  "
$ws.Range("D22").Value = "376.13"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").Value = "73.96"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "0.546"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -9.19%  "
$ws.Range("E28").Value = "  -7.21%  "
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "7.27"
$ws.Range("E30").Value = "  -6.31%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "2.04"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("D34").Value = "22.57"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -6.72%  "
$ws.Range("D36").Value = "5.25"
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("D37").Value = "166.58"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  -7.81%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.75"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "26.75"
$ws.Range("E40").Value = "  -16.02%  "
$ws.Range("D41").Value = "3.323.80"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "0.0744"
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").Value = "41.93"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("D44").Value = "0.751"
$ws.Range("E44").Value = "  -4.77%  "
$ws.Range("D45").Value = "4.21"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("E46").Value = "  -6.50%  "
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("D48").Value = "2.357.84"
$ws.Range("E48").Value = "  -7.87%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "6.54"
$ws.Range("E50").Value = "  -6.21%  "
$ws.Range("D51").Value = "0.0257"
$ws.Range("E51").Value = "  -3.98%  "

# Remove the temporary text-number-format so cells end up with no explicit style,
# same as the untouched cells around them.
$textRange.ClearFormats()
